$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Row 2 (class "1")
$ws.Range("B2").Value = 0.7711803041274439
$ws.Range("C2").Value = 0.9366754617414248
$ws.Range("D2").Value = 0.8459094519459889
$ws.Range("E2").Value = 1137

# Row 3 (class "2")
$ws.Range("B3").Value = 0.9479495268138801
$ws.Range("C3").Value = 0.9375975039001561
$ws.Range("D3").Value = 0.9427450980392157
$ws.Range("E3").Value = 641

# Row 4 (class "3")
$ws.Range("B4").Value = 0.8333333333333334
$ws.Range("C4").Value = 0.8250620347394541
$ws.Range("D4").Value = 0.8291770573566085
$ws.Range("E4").Value = 806

# Row 5 (class "4")
$ws.Range("B5").Value = 0.9448818897637795
$ws.Range("C5").Value = 0.3370786516853932
$ws.Range("D5").Value = 0.4968944099378881
$ws.Range("E5").Value = 356

# Row 6 (accuracy)
$ws.Range("B6").Value = 0.8336734693877551
$ws.Range("C6").Value = 0.8336734693877551
$ws.Range("D6").Value = 0.8336734693877551
$ws.Range("E6").Value = 0.8336734693877551

# Row 7 (macro avg)
$ws.Range("B7").Value = 0.8743362635096092
$ws.Range("C7").Value = 0.759103413016607
$ws.Range("D7").Value = 0.7786815043199253

# Row 8 (weighted avg)
$ws.Range("B8").Value = 0.8477932897629841
$ws.Range("C8").Value = 0.8336734693877551
$ws.Range("D8").Value = 0.8201733921336876
